# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp in the title row (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Abril de 2020 a las 10:52"

# --- Row 18 (Suiza): only "Casos criticos" (F) changes ---
$ws.Range("F18").Value = 204

# --- Rows 25/26: Arabia Saudita / Austria swap places with updated data ---
# Row 25 becomes Austria with new figures
$ws.Range("A25").Value = "Austria"
$ws.Range("B25").Value = 15148
$ws.Range("C25").Value = 77
$ws.Range("D25").Value = 12103
$ws.Range("E25").Value = 2509
$ws.Range("F25").Value = 148
$ws.Range("G25").Value = 6
$ws.Range("H25").Value = 536

# Row 26 becomes Arabia Saudita, keeping its previous (unchanged) figures
$ws.Range("A26").Value = "Arabia Saudita"
$ws.Range("B26").Value = 15102
$ws.Range("C26").Value = 0
$ws.Range("D26").Value = 2049
$ws.Range("E26").Value = 12926
$ws.Range("F26").Value = 93
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 127

# --- Rows 50/51: Colombia / Banglades swap places with updated data ---
# Row 50 becomes Banglades with new figures
$ws.Range("A50").Value = "Banglades"
$ws.Range("B50").Value = 4998
$ws.Range("C50").Value = 309
$ws.Range("D50").Value = 112
$ws.Range("E50").Value = 4746
$ws.Range("F50").Value = 1
$ws.Range("G50").Value = 9
$ws.Range("H50").Value = 140

# Row 51 becomes Colombia, keeping its previous (unchanged) figures
$ws.Range("A51").Value = "Colombia"
$ws.Range("B51").Value = 4881
$ws.Range("C51").Value = 0
$ws.Range("D51").Value = 1003
$ws.Range("E51").Value = 3653
$ws.Range("F51").Value = 117
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 225

# --- Row 59 (Moldavia): update activos/recuperados/muertes-hoy/muertes ---
$ws.Range("D59").Value = 825
$ws.Range("E59").Value = 2198
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 87

# --- Rows 113/114: Reunion / Consejo Danes para los Refugiados swap places with updated data ---
# Row 113 becomes Consejo Danes para los Refugiados with new figures
$ws.Range("A113").Value = "Consejo Danes para los Refugiados"
$ws.Range("B113").Value = 416
$ws.Range("C113").Value = 22
$ws.Range("D113").Value = 49
$ws.Range("E113").Value = 339
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 3
$ws.Range("H113").Value = 28

# Row 114 becomes Reunion, keeping its previous (unchanged) figures
$ws.Range("A114").Value = "Reunion"
$ws.Range("B114").Value = 412
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 300
$ws.Range("E114").Value = 112
$ws.Range("F114").Value = 2
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 0

$wb.Save()
